$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New adj_close (column G) text values for rows 2..124 (index 0 -> row 2).
# These are stored as plain text in the sheet (shared strings), mirroring
# the source data's string-typed "adj_close" column.
$adjClose = @("2.16782140731812","2.16732406616211","2.15787720680237","2.16782140731812","2.08826804161072","2.1668267250061","2.11810064315796","2.15787720680237","2.12804484367371","2.11810064315796","2.11810064315796","2.11810064315796","2.09821224212646","2.0683798789978","2.07832407951355","2.08826804161072","2.08826804161072","2.08826804161072","2.07832407951355","2.03854751586914","1.98882687091827","2.02860331535339","1.99877095222473","2.00871515274048","2.02860331535339","2.03854751586914","1.99877095222473","1.99877095222473","1.92916214466095","2.03854751586914","1.98882687091827","1.98882687091827","1.98882687091827","1.98882687091827","2.00871515274048","2.00871515274048","2.00871515274048","1.94407832622528","1.90927374362946","1.88441348075867","1.79988825321198","1.82972073554993","1.90927374362946","1.90927374362946","1.84960901737213","1.97391068935394","1.97391068935394","1.88441348075867","1.88441348075867","1.88441348075867","1.88441348075867","1.88441348075867","1.83966493606567","1.79988825321198","1.85955309867859","1.81977665424347","1.78994417190552","1.81977665424347","1.85955309867859","1.89435756206512","1.92916214466095","1.88938546180725","1.88938546180725","1.91921782493591","1.90927374362946","1.86452519893646","1.83966493606567","1.80983245372772","1.7949161529541","1.7949161529541","1.7949161529541","1.77999997138977","1.64078211784363","1.5612291097641","1.52145254611969","1.51150834560394","1.51150834560394","1.51150834560394","1.51150834560394","1.51150834560394","1.51150834560394","1.51150834560394","1.49162018299103","1.49162018299103","1.54134082794189","1.55128490924835","1.61094975471497","1.75016760826111","1.86949729919434","1.7949161529541","1.86949729919434","1.88441348075867","1.77502799034119","1.75513970851898","1.71039116382599","1.66067039966583","1.63581013679504","1.63581013679504","1.63581013679504","1.59106147289276","1.59106147289276","1.59106147289276","1.59106147289276","1.54134082794189","1.49162018299103","1.49162018299103","1.49162018299103","1.51648044586182","1.54134082794189","1.50156426429749","1.50156426429749","1.51648044586182","1.52145254611969","1.52145254611969","1.52145254611969","1.53636872768402","1.67061448097229","1.54631280899048","1.57117331027985","1.6258659362793","1.65072631835938","1.77999997138977","1.77499997615814")

# Row 123's date (column A) moved from 45464.631712963 to 45464.2916666667.
$ws.Range("A123").Value = 45464.2916666667

# Rewrite adj_close for every existing data row (2..123) as TEXT (not a
# number), matching the workbook's convention of storing this column as
# shared strings. NumberFormat "@" forces text interpretation of a
# numeric-looking string; Style is reset back to "Normal" right after so
# the cell keeps its original (default) style index.
for ($row = 2; $row -le 123; $row++) {
    $cell = $ws.Cells.Item($row, 7)
    $cell.NumberFormat = "@"
    $cell.Value = $adjClose[$row - 2]
    $cell.Style = "Normal"
}

# Append the new trading-day row (124), copying A123's date number format
# onto A124 before writing its value so it displays/stores like the other
# date cells.
$ws.Range("A123").Copy() | Out-Null
$ws.Range("A124").PasteSpecial(-4122) | Out-Null

$ws.Range("A124").Value = 45467.5717592593
$ws.Range("B124").Value = 14000
$ws.Range("C124").Value = 1.81500005722046
$ws.Range("D124").Value = 1.77499997615814
$ws.Range("E124").Value = 1.68499994277954
$ws.Range("F124").Value = 1.77499997615814

$g124 = $ws.Cells.Item(124, 7)
$g124.NumberFormat = "@"
$g124.Value = $adjClose[122]
$g124.Style = "Normal"

$ws.Range("H124").Value = "SMN.MI"

Write-Output "edit complete"
